{"js": "// The \"Main scenario\" bullet list (in the use-case table) had a stray,\n// essentially-empty bullet \"Customer inputs \" (containing only the\n// leftover \"_GoBack\" bookmark) immediately before the bullet describing\n// the system's response. This edit removes that stray bullet and fixes\n// a typo (\"constrains\" -> \"constraints\") in the following bullet.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nconst STRAY_TEXT = \"Customer inputs \";\nconst OLD_TEXT = \"System shows a list of all available flights with those constrains\";\nconst NEW_TEXT = \"System shows a list of all available flights with those constraints\";\n\nlet strayParagraph = null;\nlet targetParagraph = null;\n\nfor (const p of paragraphs.items) {\n  if (p.text === STRAY_TEXT) {\n    strayParagraph = p;\n  } else if (p.text === OLD_TEXT) {\n    targetParagraph = p;\n  }\n}\n\n// Fix the typo first (while the paragraph collection is still in its\n// original arrangement), then drop the stray bullet paragraph entirely\n// (this also removes the now-orphaned bookmarkStart/bookmarkEnd pair\n// that lived inside it).\nif (targetParagraph) {\n  targetParagraph.insertText(NEW_TEXT, \"Replace\");\n}\nif (strayParagraph) {\n  strayParagraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# The \"Main scenario\" bullet list (in the use-case table) had a stray,\n# essentially-empty bullet \"Customer inputs \" (containing only the\n# leftover \"_GoBack\" bookmark) immediately before the bullet describing\n# the system's response. This edit removes that stray bullet and fixes\n# a typo (\"constrains\" -> \"constraints\") in the following bullet.\n\n$d = $word.ActiveDocument\n\n$strayText = \"Customer inputs \"\n$oldText = \"System shows a list of all available flights with those constrains\"\n$newText = \"System shows a list of all available flights with those constraints\"\n\n$strayParagraph = $null\n$targetParagraph = $null\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t -eq $strayText) {\n        $strayParagraph = $p\n    } elseif ($t -eq $oldText) {\n        $targetParagraph = $p\n    }\n}\n\n# Fix the typo first (while the paragraph collection is still in its\n# original arrangement), then drop the stray bullet paragraph entirely\n# (this also removes the now-orphaned bookmarkStart/bookmarkEnd pair\n# that lived inside it). Assigning .Range.Text (without a trailing `r`)\n# replaces only the paragraph's own content, leaving its paragraph mark\n# untouched.\nif ($targetParagraph -ne $null) {\n    $targetParagraph.Range.Text = $newText\n}\nif ($strayParagraph -ne $null) {\n    $strayParagraph.Range.Delete()\n}\n"}
